$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 46022.01041666666, 903.744),
    @(3, 46022.02083333334, 936.202),
    @(4, 46022.03125, 964.2329999999999),
    @(5, 46022.04166666666, 1005.484),
    @(6, 46022.05208333334, 1135.639),
    @(7, 46022.0625, 1158.652),
    @(8, 46022.07291666666, 1182.285),
    @(9, 46022.08333333334, 1211.832),
    @(10, 46022.09375, 1332.306),
    @(11, 46022.10416666666, 1357.809),
    @(12, 46022.11458333334, 1378.94),
    @(13, 46022.125, 1394.862),
    @(14, 46022.13541666666, 1479.896),
    @(15, 46022.14583333334, 1476.957),
    @(16, 46022.15625, 1474.458),
    @(17, 46022.16666666666, 1458.317),
    @(18, 46022.17708333334, 1469.855),
    @(19, 46022.1875, 1462.407),
    @(20, 46022.19791666666, 1463.447),
    @(21, 46022.20833333334, 1464.583),
    @(22, 46022.21875, 1466.976),
    @(23, 46022.22916666666, 1468.8),
    @(24, 46022.23958333334, 1473.553),
    @(25, 46022.25, 1492.896),
    @(26, 46022.26041666666, 1503.965),
    @(27, 46022.27083333334, 1510.492),
    @(28, 46022.28125, 1514.385),
    @(29, 46022.29166666666, 1514.651),
    @(30, 46022.30208333334, 1519.159),
    @(31, 46022.3125, 1533.212),
    @(32, 46022.32291666666, 1543.169),
    @(33, 46022.33333333334, 1546.354),
    @(34, 46022.34375, 1590.506),
    @(35, 46022.35416666666, 1575.339),
    @(36, 46022.36458333334, 1566.141),
    @(37, 46022.375, 1551.08),
    @(38, 46022.38541666666, 1542.78),
    @(39, 46022.39583333334, 1538.739),
    @(40, 46022.40625, 1532.367),
    @(41, 46022.41666666666, 1527.225),
    @(42, 46022.42708333334, 1542.189),
    @(43, 46022.4375, 1549.313),
    @(44, 46022.44791666666, 1558.065),
    @(45, 46022.45833333334, 1566.065),
    @(46, 46022.46875, 1593.893),
    @(47, 46022.47916666666, 1603.318),
    @(48, 46022.48958333334, 1609.563),
    @(49, 46022.5, 1614.984),
    @(50, 46022.51041666666, 1631.337),
    @(51, 46022.52083333334, 1629.503),
    @(52, 46022.53125, 1628.493),
    @(53, 46022.54166666666, 1627.616),
    @(54, 46022.55208333334, 1589.559),
    @(55, 46022.5625, 1573.236),
    @(56, 46022.57291666666, 1568.323),
    @(57, 46022.58333333334, 1553.258),
    @(58, 46022.59375, 1398.222),
    @(59, 46022.60416666666, 1351.237),
    @(60, 46022.61458333334, 1319.67),
    @(61, 46022.625, 1288.401),
    @(62, 46022.63541666666, 1153.061),
    @(63, 46022.64583333334, 1148.212),
    @(64, 46022.65625, 1146.077),
    @(65, 46022.66666666666, 1145.868),
    @(66, 46022.67708333334, 1126.743),
    @(67, 46022.6875, 1133.628),
    @(68, 46022.69791666666, 1137.36),
    @(69, 46022.70833333334, 1141.682),
    @(70, 46022.71875, 1126.966),
    @(71, 46022.72916666666, 1123.421),
    @(72, 46022.73958333334, 1117.823),
    @(73, 46022.75, 1112.883),
    @(74, 46022.76041666666, 1083.455),
    @(75, 46022.77083333334, 1075.758),
    @(76, 46022.78125, 1067.595),
    @(77, 46022.79166666666, 1058.242),
    @(78, 46022.80208333334, 1003.46),
    @(79, 46022.8125, 992.689),
    @(80, 46022.82291666666, 982.527),
    @(81, 46022.83333333334, 972.771),
    @(82, 46022.84375, 914.629),
    @(83, 46022.85416666666, 904.546),
    @(84, 46022.86458333334, 893.795),
    @(85, 46022.875, 881.6900000000001),
    @(86, 46022.88541666666, 823.369),
    @(87, 46022.89583333334, 809.65),
    @(88, 46022.90625, 797.66),
    @(89, 46022.91666666666, 785.998),
    @(90, 46022.92708333334, 728.465),
    @(91, 46022.9375, 715.326),
    @(92, 46022.94791666666, 698.865),
    @(93, 46022.95833333334, 685.077),
    @(94, 46022.96875, $null),
    @(95, 46022.97916666666, $null),
    @(96, 46022.98958333334, $null),
    @(97, 46023, $null)
)

foreach ($row in $data) {
    $r = $row[0]
    $a = $row[1]
    $b = $row[2]
    $ws.Cells.Item($r, 1).Value = $a
    if ($null -ne $b) {
        $ws.Cells.Item($r, 2).Value = $b
    }
}
